$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Ensure the timestamp is written as plain text, not auto-converted to a date/number
$ws.Cells.Item($newRow, 1).Value = "'2025-04-29 08:44:44"
$ws.Cells.Item($newRow, 2).Value = 191

# Match the style of the previous data row (no special formatting / quote prefix)
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($lastRow, 1).Style
